$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.408.71"
Set-TextValue $ws.Range("E2") "  -0.77%  "
Set-TextValue $ws.Range("D3") "1.723.17"
Set-TextValue $ws.Range("E3") "  -0.46%  "
Set-TextValue $ws.Range("D4") "0.9998"
Set-TextValue $ws.Range("E4") "  +0.11%  "
Set-TextValue $ws.Range("D5") "243.95"
Set-TextValue $ws.Range("E5") "  -1.07%  "
Set-TextValue $ws.Range("E6") "  +0.11%  "
Set-TextValue $ws.Range("D7") "0.4932"
Set-TextValue $ws.Range("E7") "  +2.38%  "
Set-TextValue $ws.Range("D8") "0.2612"
Set-TextValue $ws.Range("E8") "  -2.56%  "
Set-TextValue $ws.Range("D9") "0.06211"
Set-TextValue $ws.Range("E9") "  +0.40%  "
Set-TextValue $ws.Range("D10") "1.721.76"
Set-TextValue $ws.Range("E10") "  -0.42%  "
Set-TextValue $ws.Range("D11") "0.06992"
Set-TextValue $ws.Range("E11") "  -2.31%  "
Set-TextValue $ws.Range("D12") "15.45"
Set-TextValue $ws.Range("E12") "  -1.30%  "
Set-TextValue $ws.Range("D13") "4.545"
Set-TextValue $ws.Range("E13") "  -0.03%  "
Set-TextValue $ws.Range("D14") "0.6000"
Set-TextValue $ws.Range("E14") "  -2.10%  "
Set-TextValue $ws.Range("D15") "77.46"
Set-TextValue $ws.Range("E15") "  +0.05%  "
Set-TextValue $ws.Range("D16") "1.0000"
Set-TextValue $ws.Range("E16") "  +0.02%  "
Set-TextValue $ws.Range("D17") "26.407.62"
Set-TextValue $ws.Range("E17") "  -0.74%  "
Set-TextValue $ws.Range("D18") "0.9999"
Set-TextValue $ws.Range("E18") "  +0.10%  "
Set-TextValue $ws.Range("D19") "0.000007194"
Set-TextValue $ws.Range("E19") "  +3.22%  "
Set-TextValue $ws.Range("D20") "11.35"
Set-TextValue $ws.Range("E20") "  -1.88%  "
Set-TextValue $ws.Range("D21") "1.943.88"
Set-TextValue $ws.Range("E21") "  -0.47%  "
Set-TextValue $ws.Range("D22") "4.471"
Set-TextValue $ws.Range("E22") "  -1.12%  "
Set-TextValue $ws.Range("D23") "8.577"
Set-TextValue $ws.Range("E23") "  -2.70%  "
Set-TextValue $ws.Range("D24") "5.158"
Set-TextValue $ws.Range("E24") "  -1.75%  "
Set-TextValue $ws.Range("E25") "  +0.47%  "
Set-TextValue $ws.Range("E26") "  -0.97%  "
Set-TextValue $ws.Range("E27") "  -0.75%  "
Set-TextValue $ws.Range("D28") "107.02"
Set-TextValue $ws.Range("E28") "  -1.40%  "
Set-TextValue $ws.Range("D29") "1.721"
Set-TextValue $ws.Range("E29") "  -3.53%  "
Set-TextValue $ws.Range("D30") "3.950"
Set-TextValue $ws.Range("E30") "  -0.73%  "
Set-TextValue $ws.Range("D31") "0.08020"
Set-TextValue $ws.Range("E31") "  -0.06%  "
Set-TextValue $ws.Range("D32") "3.679"
Set-TextValue $ws.Range("E32") "  -0.40%  "
Set-TextValue $ws.Range("E33") "  -0.38%  "
Set-TextValue $ws.Range("D34") "0.9992"
Set-TextValue $ws.Range("E34") "  +0.03%  "
Set-TextValue $ws.Range("E35") "  -0.29%  "
Set-TextValue $ws.Range("D36") "0.9995"
Set-TextValue $ws.Range("E36") "  -0.71%  "
Set-TextValue $ws.Range("D37") "0.6269"
Set-TextValue $ws.Range("E37") "  -1.26%  "
Set-TextValue $ws.Range("D38") "0.9437"
Set-TextValue $ws.Range("E38") "  +4.79%  "
Set-TextValue $ws.Range("D39") "2.392"
Set-TextValue $ws.Range("E39") "  +0.77%  "
Set-TextValue $ws.Range("D40") "1.946"
Set-TextValue $ws.Range("E40") "  -5.27%  "
Set-TextValue $ws.Range("D41") "0.9997"
Set-TextValue $ws.Range("E41") "  -0.37%  "
Set-TextValue $ws.Range("D42") "0.01482"
Set-TextValue $ws.Range("E42") "  -1.43%  "
Set-TextValue $ws.Range("D43") "99.60"
Set-TextValue $ws.Range("E43") "  -3.52%  "
Set-TextValue $ws.Range("D44") "5.285"
Set-TextValue $ws.Range("E44") "  -3.44%  "
Set-TextValue $ws.Range("D45") "0.3854"
Set-TextValue $ws.Range("E45") "  -1.46%  "
Set-TextValue $ws.Range("D46") "6.822"
Set-TextValue $ws.Range("E46") "  -5.12%  "
Set-TextValue $ws.Range("D47") "0.1169"
Set-TextValue $ws.Range("E47") "  -1.50%  "
Set-TextValue $ws.Range("D48") "0.05365"
Set-TextValue $ws.Range("E48") "  -0.35%  "
Set-TextValue $ws.Range("D49") "7.784"
Set-TextValue $ws.Range("E49") "  -1.42%  "
Set-TextValue $ws.Range("D50") "30.23"
Set-TextValue $ws.Range("E50") "  -1.58%  "
Set-TextValue $ws.Range("D51") "1.234"
Set-TextValue $ws.Range("E51") "  -1.72%  "
